$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C7").Value = "Worked with previously done control statement practical works, arrays: int, string, sorting: asc, dec"
